$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($findText, $true, $false, $false, $false, $false, `
                  $true, 1, $false, $replaceText, 2) | Out-Null
}

# 1. Merge the split "Identify the choice..." runs back into a single run.
Replace-Text `
    "Identify the choice of solution and approach that will be used to deliver the intended project outcomes. " `
    "Identify the choice of solution and approach that will be used to deliver the intended project outcomes. "

# 2. Merge the split "Highlight the key risks..." runs back into a single run.
Replace-Text `
    "Highlight the key risks to the project, together with the likelihood of, and strategies for mitigating, each risk." `
    "Highlight the key risks to the project, together with the likelihood of, and strategies for mitigating, each risk."

# 3. Merge the split "Provide a summary of the monitoring..." runs back into a single run.
Replace-Text `
    "Provide a summary of the monitoring and reporting mechanisms that will be adopted in relation to the project and the frequency with which each will be implemented." `
    "Provide a summary of the monitoring and reporting mechanisms that will be adopted in relation to the project and the frequency with which each will be implemented."

# 4. Update the header text: "GMIT Civic Engagement" -> "Civic Engagement"
foreach ($sec in $d.Sections) {
    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            $find = $hdr.Range.Find
            $find.ClearFormatting()
            $find.Replacement.ClearFormatting()
            $find.Execute("GMIT Civic Engagement", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "Civic Engagement", 2) | Out-Null
        }
    }
}
